# hw2_requirements.xlsx update: refresh answers/notes per the README & testcases rewrite.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 6 - "2. You should cache appropriate responses"
#   answer changed from "implement, but has bug" to "yes"
#   note changed to describe the caching policy
$ws.Range("B6").Value = "yes"
$ws.Range("C6").Value = "We cache response that are 200 OK and not private, not no-store"

# Row 7 - "2.1 You should follow rules of expiration time" - new extra note
$ws.Range("C7").Value = "Please see README for details"

# Row 8 - "2.2 You should follow rules of re-validation" - new extra note
$ws.Range("C8").Value = "Please see README for details"

# Row 11 - "3.1 Do you spawn a thread/process to handle a request?" - reworded note
$ws.Range("C11").Value = "Every time we get a request from brower, we create a new thread. And we use mutex lock when writing to log file and editting cache"

# Row 12 - "3.2 ... synchronization for your cache" - note removed
$ws.Range("C12").Value = ""

# The newly-typed notes in column C pick up the plain (non-wrapped) Arial 10 black style
# that is already used elsewhere in the sheet (e.g. the "Extra note to TA" cell).
foreach ($addr in @("C6", "C7", "C8", "C11")) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Arial"
    $cell.Font.Size = 10
    $cell.Font.ColorIndex = 1
    $cell.WrapText = $False
}

# Move the active selection to B13, matching where the user left off editing.
$ws.Range("B13").Select()
